# TODO.xlsx: remove the six "Next" / "NextNext" tasks (old rows 2-7).
# Everything below shifts up by 6 rows; the Excel Table ("Tabelle1"),
# conditional formatting and shared-string table are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 6 obsolete task rows (old rows 2-7). This shrinks the backing
# Excel Table automatically (it spans A1:E29 -> A1:E23) and shifts every
# row below up by 6.
$ws.Range("A2:E7").EntireRow.Delete()

# Conditional formatting rules reference fixed single cells (not part of
# the deleted block) that need to be re-pointed to their new row numbers.
$fcs = $ws.Cells.FormatConditions

$fc = $fcs.Item(2)
$fc.ModifyAppliesToRange($ws.Range("B11"))
$fc.Formula1 = 'NOT(ISERROR(SEARCH("Unclear",B11)))'

$fc = $fcs.Item(3)
$fc.ModifyAppliesToRange($ws.Range("E20"))
$fc.Formula1 = 'NOT(ISERROR(SEARCH("Unclear",E20)))'

$fc = $fcs.Item(4)
$fc.ModifyAppliesToRange($ws.Range("E23"))
$fc.Formula1 = 'NOT(ISERROR(SEARCH("Unclear",E23)))'

$fc = $fcs.Item(5)
$fc.ModifyAppliesToRange($ws.Range("B9"))
$fc.Formula1 = 'NOT(ISERROR(SEARCH("Unclear",B9)))'
